{"js": "// The worksheet table has 20 rows x 5 columns; only every 4th row (0, 4, 8,\n// 12, 16) holds the math-problem text, the rows in between are blank\n// \"work space\" rows. Update the 25 populated cells to their new values,\n// addressing each cell by its (row, column) position rather than by its old\n// text, since several old/new values collide or repeat across cells.\nconst newValues = [\n  [0, [\"21\u00f75=4, 1\", \"76\u00f73=25, 1\", \"42\u00f78=5, 2\", \"91\u00f77=13, 0\", \"21\u00f75=4, 1\"]],\n  [4, [\"60\u00f78=7, 4\", \"26\u00f78=3, 2\", \"15\u00f75=3, 0\", \"21\u00f73=7, 0\", \"67\u00f73=22, 1\"]],\n  [8, [\"40\u00f77=5, 5\", \"51\u00f75=10, 1\", \"78\u00f77=11, 1\", \"76\u00f79=8, 4\", \"86\u00f77=12, 2\"]],\n  [12, [\"67\u00f76=11, 1\", \"42\u00f76=7, 0\", \"65\u00f74=16, 1\", \"41\u00f76=6, 5\", \"28\u00f77=4, 0\"]],\n  [16, [\"10\u00f75=2, 0\", \"55\u00f78=6, 7\", \"89\u00f79=9, 8\", \"57\u00f76=9, 3\", \"94\u00f72=47, 0\"]],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nfor (const [rowIndex, rowValues] of newValues) {\n  for (let colIndex = 0; colIndex < rowValues.length; colIndex++) {\n    const cell = table.getCell(rowIndex, colIndex);\n    cell.value = rowValues[colIndex];\n  }\n}\nawait context.sync();\n", "ps1": "# The worksheet table has 20 rows x 5 columns; only every 4th row (1, 5, 9,\n# 13, 17 in Word's 1-based Cell() indexing) holds the math-problem text, the\n# rows in between are blank \"work space\" rows. Update the 25 populated cells\n# to their new values, addressing each cell by its (row, column) position\n# rather than by its old text, since several old/new values collide or\n# repeat across cells.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"21\u00f75=4, 1\", \"76\u00f73=25, 1\", \"42\u00f78=5, 2\", \"91\u00f77=13, 0\", \"21\u00f75=4, 1\")\n    5  = @(\"60\u00f78=7, 4\", \"26\u00f78=3, 2\", \"15\u00f75=3, 0\", \"21\u00f73=7, 0\", \"67\u00f73=22, 1\")\n    9  = @(\"40\u00f77=5, 5\", \"51\u00f75=10, 1\", \"78\u00f77=11, 1\", \"76\u00f79=8, 4\", \"86\u00f77=12, 2\")\n    13 = @(\"67\u00f76=11, 1\", \"42\u00f76=7, 0\", \"65\u00f74=16, 1\", \"41\u00f76=6, 5\", \"28\u00f77=4, 0\")\n    17 = @(\"10\u00f75=2, 0\", \"55\u00f78=6, 7\", \"89\u00f79=9, 8\", \"57\u00f76=9, 3\", \"94\u00f72=47, 0\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $rowValues = $newValues[$rowIndex]\n    for ($col = 1; $col -le $rowValues.Length; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
